$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$rowData = @(
    "8Y1cGcHD", "14/11/2024", "19:30", "URUGUAY - PRIMERA DIVISION", "Wanderers", "Penarol",
    6.25, 4.33, 1.5, 6.5, 2.25, 2.05, 1.05, 11, 1.29, 3.5, 1.95, 1.9, 1.4, 2.75, 2, 1.73,
    15, 34, 19, 67, 51, 51, 10, 8, 21, 67, 6.5, 6.5, 8.5, 10, 13, 29, 1250, 8, 34, 41,
    151, 151, 351, 2.75, 9, 67, 3.4, 7.5, 21, 23, 51, 151, 51, 51
)

for ($i = 0; $i -lt $rowData.Length; $i++) {
    $col = $i + 1
    $ws.Cells.Item(4, $col).Value = $rowData[$i]
}
